# Update patient admission/discharge form (RI10.10.2.56) with new patient
# data - para agregar codigo de barras
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = "ZZ1"

$ws.Range($scratch).Value = '="IBOY"'
$ws.Range($scratch).Copy()
$ws.Range("A6").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="SOLARES"'
$ws.Range($scratch).Copy()
$ws.Range("C6").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="ISABEL"'
$ws.Range($scratch).Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range($scratch).Value = '=""'
$ws.Range($scratch).Copy()
$ws.Range("G6").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="2016-26508/201766161"'
$ws.Range($scratch).Copy()
$ws.Range("I6").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="ALDEA PUERTA PARADA"'
$ws.Range($scratch).Copy()
$ws.Range("A8").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="SANTA CATARINA PINULA"'
$ws.Range($scratch).Copy()
$ws.Range("D8").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="GUATEMALA CALLEJON LA LUZ"'
$ws.Range($scratch).Copy()
$ws.Range("F8").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="GUATEMALA"'
$ws.Range($scratch).Copy()
$ws.Range("H8").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="54062245"'
$ws.Range($scratch).Copy()
$ws.Range("J8").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="1945-11-05"'
$ws.Range($scratch).Copy()
$ws.Range("A12").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="72"'
$ws.Range($scratch).Copy()
$ws.Range("F12").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="Femenino"'
$ws.Range($scratch).Copy()
$ws.Range("J12").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="Casado"'
$ws.Range($scratch).Copy()
$ws.Range("A14").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="AMA  DE CASA"'
$ws.Range($scratch).Copy()
$ws.Range("D14").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="GUATEMALTECA"'
$ws.Range($scratch).Copy()
$ws.Range("F14").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="1683101550102"'
$ws.Range($scratch).Copy()
$ws.Range("H14").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="SANTAANA ARIZANDIETA"'
$ws.Range($scratch).Copy()
$ws.Range("A16").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="FRANCISCO IBOY PINEDA"'
$ws.Range($scratch).Copy()
$ws.Range("A18").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="ENGRACIA SOLARES"'
$ws.Range($scratch).Copy()
$ws.Range("F18").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="NANCY ARIZANDIETA"'
$ws.Range($scratch).Copy()
$ws.Range("A20").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="HIJA"'
$ws.Range($scratch).Copy()
$ws.Range("F20").PasteSpecial(-4163)

$ws.Range($scratch).Value = '=""'
$ws.Range($scratch).Copy()
$ws.Range("H20").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="42206330"'
$ws.Range($scratch).Copy()
$ws.Range("J20").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="20/11/2017"'
$ws.Range($scratch).Copy()
$ws.Range("A24").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="14:30:22"'
$ws.Range($scratch).Copy()
$ws.Range("C24").PasteSpecial(-4163)

$ws.Range($scratch).Value = '="INGRESO A HEMATO-ONCO"'
$ws.Range($scratch).Copy()
$ws.Range("D24").PasteSpecial(-4163)

$ws.Range($scratch).ClearContents()
$excel.CutCopyMode = $false
$wb.Save()
